$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "sfasfasdfsa"
$ws.Range("C6").Value = " sfsa"

$ws.Range("A5").Select()
